$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 14.47581233333333
$ws.Range("H2").Value = 43.427437
$ws.Range("I2").Value = 0.2264097534340379
$ws.Range("J2").Value = 0.2264097534340379
$ws.Range("M2").Value = 62.58874
$ws.Range("N2").Value = 187.76622
$ws.Range("O2").Value = 0.4331197020873656
$ws.Range("P2").Value = 0.4331197020873655
$ws.Range("Q2").Value = 906.0228544197932
$ws.Range("R2").Value = 8154.205689778139
$ws.Range("S2").Value = 0.09806252495702436
$ws.Range("T2").Value = 0.09806252495702436

$ws.Range("G3").Value = 14.47581233333333
$ws.Range("H3").Value = 43.427437
$ws.Range("I3").Value = 0.2264097534340379
$ws.Range("J3").Value = 0.2264097534340379
$ws.Range("O3").Value = 0.03956530671562308
$ws.Range("P3").Value = 0.03956530671562306
$ws.Range("Q3").Value = 82.76481525481987
$ws.Range("R3").Value = 744.8833372933789
$ws.Range("S3").Value = 0.008957971338026303
$ws.Range("T3").Value = 0.0089579713380263

$ws.Range("G4").Value = 14.47581233333333
$ws.Range("H4").Value = 43.427437
$ws.Range("I4").Value = 0.2264097534340379
$ws.Range("J4").Value = 0.2264097534340379
$ws.Range("M4").Value = 16.124321
$ws.Range("N4").Value = 48.372963
$ws.Range("O4").Value = 0.111581749494894
$ws.Range("P4").Value = 0.111581749494894
$ws.Range("Q4").Value = 233.4126447984256
$ws.Range("R4").Value = 2100.713803185831
$ws.Range("S4").Value = 0.02526319639087753
$ws.Range("T4").Value = 0.02526319639087753

$ws.Range("G5").Value = 14.47581233333333
$ws.Range("H5").Value = 43.427437
$ws.Range("I5").Value = 0.2264097534340379
$ws.Range("J5").Value = 0.2264097534340379
$ws.Range("M5").Value = 11.37633566666667
$ws.Range("N5").Value = 34.129007
$ws.Range("O5").Value = 0.07872526455705194
$ws.Range("P5").Value = 0.07872526455705191
$ws.Range("Q5").Value = 164.6817001516732
$ws.Range("R5").Value = 1482.135301365059
$ws.Range("S5").Value = 0.01782416773739153
$ws.Range("T5").Value = 0.01782416773739152

$ws.Range("G6").Value = 14.47581233333333
$ws.Range("H6").Value = 43.427437
$ws.Range("I6").Value = 0.2264097534340379
$ws.Range("J6").Value = 0.2264097534340379
$ws.Range("M6").Value = 36.54706633333333
$ws.Range("N6").Value = 109.641199
$ws.Range("O6").Value = 0.2529089814311731
$ws.Range("P6").Value = 0.2529089814311731
$ws.Range("Q6").Value = 529.048473575218
$ws.Range("R6").Value = 4761.436262176962
$ws.Range("S6").Value = 0.05726106012708555
$ws.Range("T6").Value = 0.05726106012708556

$ws.Range("G7").Value = 14.47581233333333
$ws.Range("H7").Value = 43.427437
$ws.Range("I7").Value = 0.2264097534340379
$ws.Range("J7").Value = 0.2264097534340379
$ws.Range("M7").Value = 12.15287633333333
$ws.Range("N7").Value = 36.458629
$ws.Range("O7").Value = 0.08409899571389245
$ws.Range("P7").Value = 0.08409899571389244
$ws.Range("Q7").Value = 175.9227571115414
$ws.Range("R7").Value = 1583.304814003873
$ws.Range("S7").Value = 0.0190408328836326
$ws.Range("T7").Value = 0.01904083288363259

$ws.Range("I8").Value = 0.1838014431832978
$ws.Range("J8").Value = 0.1838014431832978
$ws.Range("M8").Value = 62.58874
$ws.Range("N8").Value = 187.76622
$ws.Range("O8").Value = 0.4331197020873656
$ws.Range("P8").Value = 0.4331197020873655
$ws.Range("Q8").Value = 735.5173779999067
$ws.Range("R8").Value = 6619.65640199916
$ws.Range("S8").Value = 0.07960802631477777
$ws.Range("T8").Value = 0.07960802631477777

$ws.Range("I9").Value = 0.1838014431832978
$ws.Range("J9").Value = 0.1838014431832978
$ws.Range("O9").Value = 0.03956530671562308
$ws.Range("P9").Value = 0.03956530671562306
$ws.Range("S9").Value = 0.007272160474321345
$ws.Range("T9").Value = 0.007272160474321342

$ws.Range("I10").Value = 0.1838014431832978
$ws.Range("J10").Value = 0.1838014431832978
$ws.Range("M10").Value = 16.124321
$ws.Range("N10").Value = 48.372963
$ws.Range("O10").Value = 0.111581749494894
$ws.Range("P10").Value = 0.111581749494894
$ws.Range("Q10").Value = 189.4864524185793
$ws.Range("R10").Value = 1705.378071767214
$ws.Range("S10").Value = 0.02050888659007872
$ws.Range("T10").Value = 0.02050888659007872

$ws.Range("I11").Value = 0.1838014431832978
$ws.Range("J11").Value = 0.1838014431832978
$ws.Range("M11").Value = 11.37633566666667
$ws.Range("N11").Value = 34.129007
$ws.Range("O11").Value = 0.07872526455705194
$ws.Range("P11").Value = 0.07872526455705191
$ws.Range("Q11").Value = 133.6900627939385
$ws.Range("R11").Value = 1203.210565145446
$ws.Range("S11").Value = 0.01446981724057307
$ws.Range("T11").Value = 0.01446981724057306

$ws.Range("I12").Value = 0.1838014431832978
$ws.Range("J12").Value = 0.1838014431832978
$ws.Range("M12").Value = 36.54706633333333
$ws.Range("N12").Value = 109.641199
$ws.Range("O12").Value = 0.2529089814311731
$ws.Range("P12").Value = 0.2529089814311731
$ws.Range("Q12").Value = 429.4862367109802
$ws.Range("R12").Value = 3865.376130398822
$ws.Range("S12").Value = 0.04648503578106747
$ws.Range("T12").Value = 0.04648503578106747

$ws.Range("I13").Value = 0.1838014431832978
$ws.Range("J13").Value = 0.1838014431832978
$ws.Range("M13").Value = 12.15287633333333
$ws.Range("N13").Value = 36.458629
$ws.Range("O13").Value = 0.08409899571389245
$ws.Range("P13").Value = 0.08409899571389244
$ws.Range("Q13").Value = 142.8156523977069
$ws.Range("R13").Value = 1285.340871579362
$ws.Range("S13").Value = 0.01545751678247941
$ws.Range("T13").Value = 0.01545751678247941

$ws.Range("G14").Value = 5.002325333333333
$ws.Range("H14").Value = 15.006976
$ws.Range("I14").Value = 0.07823914950243377
$ws.Range("J14").Value = 0.07823914950243377
$ws.Range("M14").Value = 62.58874
$ws.Range("N14").Value = 187.76622
$ws.Range("O14").Value = 0.4331197020873656
$ws.Range("P14").Value = 0.4331197020873655
$ws.Range("Q14").Value = 313.0892396834133
$ws.Range("R14").Value = 2817.80315715072
$ws.Range("S14").Value = 0.03388691712406297
$ws.Range("T14").Value = 0.03388691712406296

$ws.Range("G15").Value = 5.002325333333333
$ws.Range("H15").Value = 15.006976
$ws.Range("I15").Value = 0.07823914950243377
$ws.Range("J15").Value = 0.07823914950243377
$ws.Range("O15").Value = 0.03956530671562308
$ws.Range("P15").Value = 0.03956530671562306
$ws.Range("Q15").Value = 28.60057332357689
$ws.Range("R15").Value = 257.4051599121919
$ws.Range("S15").Value = 0.003095555947233281
$ws.Range("T15").Value = 0.003095555947233279

$ws.Range("G16").Value = 5.002325333333333
$ws.Range("H16").Value = 15.006976
$ws.Range("I16").Value = 0.07823914950243377
$ws.Range("J16").Value = 0.07823914950243377
$ws.Range("M16").Value = 16.124321
$ws.Range("N16").Value = 48.372963
$ws.Range("O16").Value = 0.111581749494894
$ws.Range("P16").Value = 0.111581749494894
$ws.Range("Q16").Value = 80.65909942109866
$ws.Range("R16").Value = 725.931894789888
$ws.Range("S16").Value = 0.008730061180474124
$ws.Range("T16").Value = 0.008730061180474124

$ws.Range("G17").Value = 5.002325333333333
$ws.Range("H17").Value = 15.006976
$ws.Range("I17").Value = 0.07823914950243377
$ws.Range("J17").Value = 0.07823914950243377
$ws.Range("M17").Value = 11.37633566666667
$ws.Range("N17").Value = 34.129007
$ws.Range("O17").Value = 0.07872526455705194
$ws.Range("P17").Value = 0.07872526455705191
$ws.Range("Q17").Value = 56.90813210587023
$ws.Range("R17").Value = 512.173188952832
$ws.Range("S17").Value = 0.006159397743297837
$ws.Range("T17").Value = 0.006159397743297834

$ws.Range("G18").Value = 5.002325333333333
$ws.Range("H18").Value = 15.006976
$ws.Range("I18").Value = 0.07823914950243377
$ws.Range("J18").Value = 0.07823914950243377
$ws.Range("M18").Value = 36.54706633333333
$ws.Range("N18").Value = 109.641199
$ws.Range("O18").Value = 0.2529089814311731
$ws.Range("P18").Value = 0.2529089814311731
$ws.Range("Q18").Value = 182.8203157782471
$ws.Range("R18").Value = 1645.382842004224
$ws.Range("S18").Value = 0.0197873836087018
$ws.Range("T18").Value = 0.0197873836087018

$ws.Range("G19").Value = 5.002325333333333
$ws.Range("H19").Value = 15.006976
$ws.Range("I19").Value = 0.07823914950243377
$ws.Range("J19").Value = 0.07823914950243377
$ws.Range("M19").Value = 12.15287633333333
$ws.Range("N19").Value = 36.458629
$ws.Range("O19").Value = 0.08409899571389245
$ws.Range("P19").Value = 0.08409899571389244
$ws.Range("Q19").Value = 60.79264115510045
$ws.Range("R19").Value = 547.133770395904
$ws.Range("S19").Value = 0.006579833898663769
$ws.Range("T19").Value = 0.006579833898663768

$ws.Range("G20").Value = 7.937383666666666
$ws.Range("H20").Value = 23.812151
$ws.Range("I20").Value = 0.1241450937259797
$ws.Range("J20").Value = 0.1241450937259797
$ws.Range("M20").Value = 62.58874
$ws.Range("N20").Value = 187.76622
$ws.Range("O20").Value = 0.4331197020873656
$ws.Range("P20").Value = 0.4331197020873655
$ws.Range("Q20").Value = 496.7908425932467
$ws.Range("R20").Value = 4471.11758333922
$ws.Range("S20").Value = 0.0537696860102044
$ws.Range("T20").Value = 0.0537696860102044

$ws.Range("G21").Value = 7.937383666666666
$ws.Range("H21").Value = 23.812151
$ws.Range("I21").Value = 0.1241450937259797
$ws.Range("J21").Value = 0.1241450937259797
$ws.Range("O21").Value = 0.03956530671562308
$ws.Range("P21").Value = 0.03956530671562306
$ws.Range("Q21").Value = 45.38163922349077
$ws.Range("R21").Value = 408.4347530114169
$ws.Range("S21").Value = 0.00491183871050816
$ws.Range("T21").Value = 0.004911838710508159

$ws.Range("G22").Value = 7.937383666666666
$ws.Range("H22").Value = 23.812151
$ws.Range("I22").Value = 0.1241450937259797
$ws.Range("J22").Value = 0.1241450937259797
$ws.Range("M22").Value = 16.124321
$ws.Range("N22").Value = 48.372963
$ws.Range("O22").Value = 0.111581749494894
$ws.Range("P22").Value = 0.111581749494894
$ws.Range("Q22").Value = 127.9849221414903
$ws.Range("R22").Value = 1151.864299273413
$ws.Range("S22").Value = 0.0138523267491524
$ws.Range("T22").Value = 0.0138523267491524

$ws.Range("G23").Value = 7.937383666666666
$ws.Range("H23").Value = 23.812151
$ws.Range("I23").Value = 0.1241450937259797
$ws.Range("J23").Value = 0.1241450937259797
$ws.Range("M23").Value = 11.37633566666667
$ws.Range("N23").Value = 34.129007
$ws.Range("O23").Value = 0.07872526455705194
$ws.Range("P23").Value = 0.07872526455705191
$ws.Range("Q23").Value = 90.29834090711745
$ws.Range("R23").Value = 812.685068164057
$ws.Range("S23").Value = 0.009773355347037759
$ws.Range("T23").Value = 0.009773355347037757

$ws.Range("G24").Value = 7.937383666666666
$ws.Range("H24").Value = 23.812151
$ws.Range("I24").Value = 0.1241450937259797
$ws.Range("J24").Value = 0.1241450937259797
$ws.Range("M24").Value = 36.54706633333333
$ws.Range("N24").Value = 109.641199
$ws.Range("O24").Value = 0.2529089814311731
$ws.Range("P24").Value = 0.2529089814311731
$ws.Range("Q24").Value = 290.0880873787832
$ws.Range("R24").Value = 2610.792786409049
$ws.Range("S24").Value = 0.03139740920391504
$ws.Range("T24").Value = 0.03139740920391504

$ws.Range("G25").Value = 7.937383666666666
$ws.Range("H25").Value = 23.812151
$ws.Range("I25").Value = 0.1241450937259797
$ws.Range("J25").Value = 0.1241450937259797
$ws.Range("M25").Value = 12.15287633333333
$ws.Range("N25").Value = 36.458629
$ws.Range("O25").Value = 0.08409899571389245
$ws.Range("P25").Value = 0.08409899571389244
$ws.Range("Q25").Value = 96.46204211121989
$ws.Range("R25").Value = 868.158379000979
$ws.Range("S25").Value = 0.01044047770516194
$ws.Range("T25").Value = 0.01044047770516194

$ws.Range("G26").Value = 5.471644666666666
$ws.Range("H26").Value = 16.414934
$ws.Range("I26").Value = 0.08557956481696133
$ws.Range("J26").Value = 0.08557956481696134
$ws.Range("M26").Value = 62.58874
$ws.Range("N26").Value = 187.76622
$ws.Range("O26").Value = 0.4331197020873656
$ws.Range("P26").Value = 0.4331197020873655
$ws.Range("Q26").Value = 342.4633454143867
$ws.Range("R26").Value = 3082.17010872948
$ws.Range("S26").Value = 0.03706619561828868
$ws.Range("T26").Value = 0.03706619561828868

$ws.Range("G27").Value = 5.471644666666666
$ws.Range("H27").Value = 16.414934
$ws.Range("I27").Value = 0.08557956481696133
$ws.Range("J27").Value = 0.08557956481696134
$ws.Range("O27").Value = 0.03956530671562308
$ws.Range("P27").Value = 0.03956530671562306
$ws.Range("Q27").Value = 31.28388580541977
$ws.Range("R27").Value = 281.554972248778
$ws.Range("S27").Value = 0.00338598173057262
$ws.Range("T27").Value = 0.00338598173057262

$ws.Range("G28").Value = 5.471644666666666
$ws.Range("H28").Value = 16.414934
$ws.Range("I28").Value = 0.08557956481696133
$ws.Range("J28").Value = 0.08557956481696134
$ws.Range("M28").Value = 16.124321
$ws.Range("N28").Value = 48.372963
$ws.Range("O28").Value = 0.111581749494894
$ws.Range("P28").Value = 0.111581749494894
$ws.Range("Q28").Value = 88.22655500327132
$ws.Range("R28").Value = 794.0389950294419
$ws.Range("S28").Value = 0.009549117563288223
$ws.Range("T28").Value = 0.009549117563288223

$ws.Range("G29").Value = 5.471644666666666
$ws.Range("H29").Value = 16.414934
$ws.Range("I29").Value = 0.08557956481696133
$ws.Range("J29").Value = 0.08557956481696134
$ws.Range("M29").Value = 11.37633566666667
$ws.Range("N29").Value = 34.129007
$ws.Range("O29").Value = 0.07872526455705194
$ws.Range("P29").Value = 0.07872526455705191
$ws.Range("Q29").Value = 62.24726637672644
$ws.Range("R29").Value = 560.225397390538
$ws.Range("S29").Value = 0.006737273880892654
$ws.Range("T29").Value = 0.006737273880892653

$ws.Range("G30").Value = 5.471644666666666
$ws.Range("H30").Value = 16.414934
$ws.Range("I30").Value = 0.08557956481696133
$ws.Range("J30").Value = 0.08557956481696134
$ws.Range("M30").Value = 36.54706633333333
$ws.Range("N30").Value = 109.641199
$ws.Range("O30").Value = 0.2529089814311731
$ws.Range("P30").Value = 0.2529089814311731
$ws.Range("Q30").Value = 199.9725605850962
$ws.Range("R30").Value = 1799.753045265866
$ws.Range("S30").Value = 0.02164384056918075
$ws.Range("T30").Value = 0.02164384056918075

$ws.Range("G31").Value = 5.471644666666666
$ws.Range("H31").Value = 16.414934
$ws.Range("I31").Value = 0.08557956481696133
$ws.Range("J31").Value = 0.08557956481696134
$ws.Range("M31").Value = 12.15287633333333
$ws.Range("N31").Value = 36.458629
$ws.Range("O31").Value = 0.08409899571389245
$ws.Range("P31").Value = 0.08409899571389244
$ws.Range("Q31").Value = 66.49622097394288
$ws.Range("R31").Value = 598.465988765486
$ws.Range("S31").Value = 0.007197155454738412
$ws.Range("T31").Value = 0.007197155454738412

$ws.Range("G32").Value = 19.29758733333334
$ws.Range("H32").Value = 57.892762
$ws.Range("I32").Value = 0.3018249953372896
$ws.Range("J32").Value = 0.3018249953372896
$ws.Range("M32").Value = 62.58874
$ws.Range("N32").Value = 187.76622
$ws.Range("O32").Value = 0.4331197020873656
$ws.Range("P32").Value = 0.4331197020873655
$ws.Range("Q32").Value = 1207.811676233294
$ws.Range("R32").Value = 10870.30508609964
$ws.Range("S32").Value = 0.1307263520630074
$ws.Range("T32").Value = 0.1307263520630074

$ws.Range("G33").Value = 19.29758733333334
$ws.Range("H33").Value = 57.892762
$ws.Range("I33").Value = 0.3018249953372896
$ws.Range("J33").Value = 0.3018249953372896
$ws.Range("O33").Value = 0.03956530671562308
$ws.Range("P33").Value = 0.03956530671562306
$ws.Range("Q33").Value = 110.3331000519616
$ws.Range("R33").Value = 992.997900467654
$ws.Range("S33").Value = 0.01194179851496137
$ws.Range("T33").Value = 0.01194179851496137

$ws.Range("G34").Value = 19.29758733333334
$ws.Range("H34").Value = 57.892762
$ws.Range("I34").Value = 0.3018249953372896
$ws.Range("J34").Value = 0.3018249953372896
$ws.Range("M34").Value = 16.124321
$ws.Range("N34").Value = 48.372963
$ws.Range("O34").Value = 0.111581749494894
$ws.Range("P34").Value = 0.111581749494894
$ws.Range("Q34").Value = 311.1604926882007
$ws.Range("R34").Value = 2800.444434193806
$ws.Range("S34").Value = 0.033678161021023
$ws.Range("T34").Value = 0.03367816102102299

$ws.Range("G35").Value = 19.29758733333334
$ws.Range("H35").Value = 57.892762
$ws.Range("I35").Value = 0.3018249953372896
$ws.Range("J35").Value = 0.3018249953372896
$ws.Range("M35").Value = 11.37633566666667
$ws.Range("N35").Value = 34.129007
$ws.Range("O35").Value = 0.07872526455705194
$ws.Range("P35").Value = 0.07872526455705191
$ws.Range("Q35").Value = 219.5358310608149
$ws.Range("R35").Value = 1975.822479547334
$ws.Range("S35").Value = 0.02376125260785909
$ws.Range("T35").Value = 0.02376125260785908

$ws.Range("G36").Value = 19.29758733333334
$ws.Range("H36").Value = 57.892762
$ws.Range("I36").Value = 0.3018249953372896
$ws.Range("J36").Value = 0.3018249953372896
$ws.Range("M36").Value = 36.54706633333333
$ws.Range("N36").Value = 109.641199
$ws.Range("O36").Value = 0.2529089814311731
$ws.Range("P36").Value = 0.2529089814311731
$ws.Range("Q36").Value = 705.2702043446266
$ws.Range("R36").Value = 6347.431839101639
$ws.Range("S36").Value = 0.07633425214122248
$ws.Range("T36").Value = 0.07633425214122248

$ws.Range("G37").Value = 19.29758733333334
$ws.Range("H37").Value = 57.892762
$ws.Range("I37").Value = 0.3018249953372896
$ws.Range("J37").Value = 0.3018249953372896
$ws.Range("M37").Value = 12.15287633333333
$ws.Range("N37").Value = 36.458629
$ws.Range("O37").Value = 0.08409899571389245
$ws.Range("P37").Value = 0.08409899571389244
$ws.Range("Q37").Value = 234.5211923936998
$ws.Range("R37").Value = 2110.690731543298
$ws.Range("S37").Value = 0.02538317898921633
$ws.Range("T37").Value = 0.02538317898921633
